$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D10: "Κείμενο άσκησης" (exercise text/question) - cleaned-up template text
$ws.Range("D10").Value = "Για τον υπολογιστή 172.16.150.10/20 να υπολογίσετε:  <br>`nΔ1. Την μάσκα δικτύου(δυαδική-δεκαδική)  <br>`nΔ2. Τη διεύθυνση δικτύου (network address)  <br> `nΔ3. Τη διεύθυνση εκπομπής (broadcast address)  <br>`nΔ4. Τον συνολικό αριθμό υπολογιστών του συγκεκριμένου δικτύου  <br>`nΔ5. Την περιοχή διευθύνσεων για υπολογιστές (από-έως) οι οποίες ανήκουν στο ίδιο δίκτυο με τον συγκεκριμένο υπολογιστή "

# E10: "Λύση άσκησης" (exercise solution) - cleaned-up long filled-in answer text
$ws.Range("E10").Value = "Για τον υπολογιστή 172.16.150.10/20 να υπολογίσετε:  <br>`nΔ1. Την μάσκα δικτύου(δυαδική-δεκαδική)<br>`n255.255.240.0 ή 11111111.11111111.1111 **0000.00000000** <br>`n`nΔ2. Τη διεύθυνση δικτύου (network address)  <br>`nΔιεύθυνση Δικτύου   <br>`n172.16.150.10(1010 **0110.00001010**) AND 255.255.240.0(1111 **00000.00000000**)=172.16.144.0(1001 **0000.00000000**)/20 <br>`n`nΔ3. Τη διεύθυνση εκπομπής (broadcast address)  <br>`nΔιεύθυνση Εκπομπής ->  172.16.159.255(1001 **1111.11111111**)  <br>`n`nΔ4. Τον συνολικό αριθμό υπολογιστών του συγκεκριμένου δικτύου  <br>`n2^12 -2 = 4094 `nΔ5. Την περιοχή διευθύνσεων για υπολογιστές (από-έως) οι οποίες ανήκουν στο ίδιο δίκτυο με τον συγκεκριμένο υπολογιστή  <br>`nΑπό 172.16.144.1 έως 172.16.159.254  <br>"

# Setting multi-line values auto-resizes the row (autofit); restore the original
# row height so row 10 stays as it was before the edit.
$ws.Rows(10).RowHeight = 14.25
